$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values look numeric, so Excel
# keeps them as literal text (matching the original inline-string data)
# instead of auto-converting to numbers.
$textCells = @("D5","D6","D7","D8","D9","D11","D12","D14","D16","D19","D20","D21","D22","D23","D24","D25","D27","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D47","D48","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values cell by cell.
$ws.Range("D2").Value = '87.813.59'
$ws.Range("E2").Value = '  +8.26%  '
$ws.Range("D3").Value = '3.322.20'
$ws.Range("E3").Value = '  +4.37%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '218.03'
$ws.Range("E5").Value = '  +4.41%  '
$ws.Range("D6").Value = '651.09'
$ws.Range("E6").Value = '  +2.77%  '
$ws.Range("D7").Value = '0.355'
$ws.Range("E7").Value = '  +23.09%  '
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").Value = '0.607'
$ws.Range("E9").Value = '  +2.53%  '
$ws.Range("D10").Value = '3.316.37'
$ws.Range("E10").Value = '  +4.13%  '
$ws.Range("D11").Value = '0.587'
$ws.Range("E11").Value = '  -0.95%  '
$ws.Range("D12").Value = '0.0000269'
$ws.Range("E12").Value = '  +1.75%  '
$ws.Range("E13").Value = '  +1.66%  '
$ws.Range("D14").Value = '35.57'
$ws.Range("E14").Value = '  +10.72%  '
$ws.Range("D15").Value = '3.924.94'
$ws.Range("E15").Value = '  +3.91%  '
$ws.Range("D16").Value = '5.48'
$ws.Range("E16").Value = '  +1.82%  '
$ws.Range("D17").Value = '87.729.26'
$ws.Range("E17").Value = '  +8.06%  '
$ws.Range("D18").Value = '3.294.35'
$ws.Range("E18").Value = '  +2.92%  '
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").Value = '14.71'
$ws.Range("E19").Value = '  +2.42%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = '9.94'
$ws.Range("E20").Value = '  +7.54%  '
$ws.Range("D21").Value = '3.13'
$ws.Range("E21").Value = '  -3.28%  '
$ws.Range("D22").Value = '456.42'
$ws.Range("E22").Value = '  +3.57%  '
$ws.Range("D23").Value = '5.54'
$ws.Range("E23").Value = '  +5.67%  '
$ws.Range("D24").Value = '5.52'
$ws.Range("E24").Value = '  +8.77%  '
$ws.Range("D25").Value = '12.61'
$ws.Range("E25").Value = '  +11.70%  '
$ws.Range("D26").Value = '3.486.52'
$ws.Range("E26").Value = '  +3.83%  '
$ws.Range("D27").Value = '78.90'
$ws.Range("E27").Value = '  +2.39%  '
$ws.Range("E28").Value = '  +0.23%  '
$ws.Range("B29").Value = 'PEPE'
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D29").Value = '0.0000126'
$ws.Range("E29").Value = '  -1.95%  '
$ws.Range("B30").Value = 'Cronos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D30").Value = '0.188'
$ws.Range("E30").Value = '  +35.50%  '
$ws.Range("D31").Value = '607.83'
$ws.Range("E31").Value = '  +6.70%  '
$ws.Range("D32").Value = '9.39'
$ws.Range("E32").Value = '  +2.20%  '
$ws.Range("D33").Value = '1.62'
$ws.Range("E33").Value = '  +7.49%  '
$ws.Range("D34").Value = '1.01'
$ws.Range("E34").Value = '  +1.09%  '
$ws.Range("D35").Value = '2.10'
$ws.Range("E35").Value = '  +2.88%  '
$ws.Range("D36").Value = '7.19'
$ws.Range("E36").Value = '  +20.78%  '
$ws.Range("D37").Value = '0.147'
$ws.Range("E37").Value = '  -4.08%  '
$ws.Range("D38").Value = '23.39'
$ws.Range("E38").Value = '  +1.59%  '
$ws.Range("D39").Value = '2.14'
$ws.Range("E39").Value = '  +5.96%  '
$ws.Range("B40").Value = 'PolygonEcosystemToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D40").Value = '0.420'
$ws.Range("E40").Value = '  +1.56%  '
$ws.Range("B41").Value = 'WhiteBITCoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D41").Value = '21.85'
$ws.Range("E41").Value = '  +5.13%  '
$ws.Range("D42").Value = '0.998'
$ws.Range("E42").Value = '  -0.22%  '
$ws.Range("D43").Value = '3.00'
$ws.Range("E43").Value = '  -3.17%  '
$ws.Range("D44").Value = '159.46'
$ws.Range("E44").Value = '  -0.42%  '
$ws.Range("D45").Value = '191.95'
$ws.Range("E45").Value = '  +1.52%  '
$ws.Range("E46").Value = '  -0.03%  '
$ws.Range("D47").Value = '1.41'
$ws.Range("E47").Value = '  +5.14%  '
$ws.Range("D48").Value = '46.20'
$ws.Range("E48").Value = '  +3.57%  '
$ws.Range("E49").Value = '  +3.12%  '
$ws.Range("D50").Value = '0.782'
$ws.Range("E50").Value = '  +0.30%  '
$ws.Range("D51").Value = '0.661'
$ws.Range("E51").Value = '  +2.73%  '
